$d = $word.ActiveDocument

# Locate the paragraph that holds the "<img id=carlJung />" marker - the
# new "styled components" block (a blank spacer + a bold-labelled pull
# quote + a trailing bold-labelled spacer) is inserted right after it.
$markerIndex = 0
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text -like "*carlJung*") {
        $markerIndex = $i
    }
}

if ($markerIndex -eq 0) {
    throw "Could not find paragraph containing the carlJung marker"
}

$marker = $d.Paragraphs.Item($markerIndex)
$insertionPoint = $d.Range($marker.Range.End, $marker.Range.End)

# Raw OOXML for the 3 new paragraphs, expressed as a pkg:package payload
# (same shape Word itself emits from Range.WordOpenXML) so every run's
# sz/b/rtl formatting lands byte-exact instead of inheriting whatever the
# insertion point happened to have. A final empty "terminator" <w:p/> is
# required: InsertXML merges the runs of the *last* paragraph in the
# fragment into the paragraph that follows the insertion point, so without
# a run-less terminator the trailing bold spacer paragraph would lose its
# own paragraph formatting by merging into the existing next paragraph.
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r><w:rPr><w:rtl w:val="0"/></w:rPr></w:r></w:p><w:p><w:pPr><w:rPr><w:b w:val="1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">&lt;quote&gt;”</w:t></w:r><w:r><w:rPr><w:sz w:val="21"/><w:szCs w:val="21"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">Lorem ipsum dolor sit amet, consectetur adipiscing elit. Fusce sagittis mollis dapibus. Duis imperdiet est odio, nec tincidunt elit vulputate ac. Sed non orci enim. Curabitur et mauris ac urna dapibus aliquam. Donec venenatis eget odio vitae suscipit. Curabitur id commodo nibh. Nulla volutpat lectus id congue ullamcorper. Sed in ipsum ut nisl venenatis tincidunt quis vel urna. Donec eu lorem urna.</w:t></w:r><w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">”&lt;/quote&gt; </w:t><w:br w:type="textWrapping"/></w:r><w:r><w:rPr><w:rtl w:val="0"/></w:rPr></w:r></w:p><w:p><w:pPr><w:rPr><w:b w:val="1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r><w:rPr><w:rtl w:val="0"/></w:rPr></w:r></w:p><w:p></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$insertionPoint.InsertXML($xml)

# The run-less terminator paragraph above becomes a genuine extra empty
# paragraph once inserted (it does not merge because it has no runs) -
# remove that scaffolding paragraph now that it has served its purpose.
# (Range.Text on an empty paragraph is just the 1-char paragraph mark.)
$strayIndex = $markerIndex + 4
$stray = $d.Paragraphs.Item($strayIndex)
if ($stray.Range.Text.Trim().Length -ne 0) {
    throw "Unexpected content in scaffolding paragraph; aborting delete"
}
$stray.Range.Delete()

Write-Host "Inserted styled quote block after paragraph $markerIndex"
